$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(100, 8).Value = 2564.9312
$ws.Cells.Item(100, 9).Value = 1778.6471
$ws.Cells.Item(100, 10).Value = 3678.8333
$ws.Cells.Item(100, 11).Value = 1778.6471
$ws.Cells.Item(100, 12).Value = 3678.8333
$ws.Cells.Item(100, 13).Value = -1237.6471
$ws.Cells.Item(100, 14).Value = -4760.8333
$ws.Cells.Item(113, 8).Value = 1802.1904
$ws.Cells.Item(113, 9).Value = 1695.875
$ws.Cells.Item(113, 11).Value = 1695.875
$ws.Cells.Item(113, 13).Value = 1558.125
$ws.Cells.Item(119, 8).Value = 802.5
$ws.Cells.Item(119, 10).Value = 802.5
$ws.Cells.Item(119, 12).Value = 2407.5
$ws.Cells.Item(119, 14).Value = -12083.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 1172.2727
$ws.Cells.Item(2, 9).Value = 952.2105
$ws.Cells.Item(2, 10).Value = 2566
$ws.Cells.Item(2, 11).Value = 952.2105
$ws.Cells.Item(2, 12).Value = 2566
$ws.Cells.Item(2, 13).Value = -839.2105
$ws.Cells.Item(2, 14).Value = -2792
$ws.Cells.Item(32, 8).Value = 2728.7
$ws.Cells.Item(32, 9).Value = 2279.281
$ws.Cells.Item(32, 10).Value = 6364.909
$ws.Cells.Item(32, 11).Value = 2279.281
$ws.Cells.Item(32, 12).Value = 6364.909
$ws.Cells.Item(32, 13).Value = -1992.281
$ws.Cells.Item(32, 14).Value = -6938.909
$ws.Cells.Item(45, 8).Value = 971.8261
$ws.Cells.Item(45, 9).Value = 953.42224
$ws.Cells.Item(45, 11).Value = 953.42224
$ws.Cells.Item(45, 13).Value = -576.42224
$ws.Cells.Item(63, 8).Value = 2291.4243
$ws.Cells.Item(63, 9).Value = 2177.5715
$ws.Cells.Item(63, 10).Value = 2490.6667
$ws.Cells.Item(63, 11).Value = 2177.5715
$ws.Cells.Item(63, 12).Value = 2490.6667
$ws.Cells.Item(63, 13).Value = -1491.5715
$ws.Cells.Item(63, 14).Value = -3862.6667
$ws.Cells.Item(64, 8).Value = 37091
$ws.Cells.Item(64, 10).Value = 37091
$ws.Cells.Item(64, 12).Value = 37091
$ws.Cells.Item(64, 14).Value = -37587
$ws.Cells.Item(66, 8).Value = 2291.4243
$ws.Cells.Item(66, 9).Value = 2177.5715
$ws.Cells.Item(66, 10).Value = 2490.6667
$ws.Cells.Item(66, 11).Value = 10887.8575
$ws.Cells.Item(66, 12).Value = 12453.3335
$ws.Cells.Item(66, 13).Value = -7455.8575
$ws.Cells.Item(66, 14).Value = -19317.3335
$ws.Cells.Item(67, 8).Value = 37091
$ws.Cells.Item(67, 10).Value = 37091
$ws.Cells.Item(67, 12).Value = 37091
$ws.Cells.Item(67, 14).Value = -38807
$ws.Cells.Item(110, 8).Value = 1087.2
$ws.Cells.Item(110, 9).Value = 1165
$ws.Cells.Item(110, 10).Value = 970.5
$ws.Cells.Item(110, 11).Value = 1165
$ws.Cells.Item(110, 12).Value = 970.5
$ws.Cells.Item(110, 13).Value = 880
$ws.Cells.Item(110, 14).Value = -5060.5
$ws.Cells.Item(116, 8).Value = 1172.2727
$ws.Cells.Item(116, 9).Value = 952.2105
$ws.Cells.Item(116, 10).Value = 2566
$ws.Cells.Item(116, 11).Value = 952.2105
$ws.Cells.Item(116, 12).Value = 2566
$ws.Cells.Item(116, 13).Value = 1341.7895
$ws.Cells.Item(116, 14).Value = -7154
$ws.Cells.Item(122, 8).Value = 1782.9412
$ws.Cells.Item(122, 9).Value = 1615.7142
$ws.Cells.Item(122, 10).Value = 1900
$ws.Cells.Item(122, 11).Value = 4847.142599999999
$ws.Cells.Item(122, 12).Value = 5700
$ws.Cells.Item(122, 13).Value = -2397.142599999999
$ws.Cells.Item(122, 14).Value = -10600

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 1172.2727
$ws.Cells.Item(3, 9).Value = 952.2105
$ws.Cells.Item(3, 10).Value = 2566
$ws.Cells.Item(3, 11).Value = 952.2105
$ws.Cells.Item(3, 12).Value = 2566
$ws.Cells.Item(3, 13).Value = -838.2105
$ws.Cells.Item(3, 14).Value = -2794
$ws.Cells.Item(98, 8).Value = 37542
$ws.Cells.Item(98, 10).Value = 37542
$ws.Cells.Item(98, 12).Value = 37542
$ws.Cells.Item(98, 14).Value = -43532
$ws.Cells.Item(107, 8).Value = 1311.4286
$ws.Cells.Item(107, 9).Value = 1198.3334
$ws.Cells.Item(107, 11).Value = 1198.3334
$ws.Cells.Item(107, 13).Value = 721.6666

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(99, 8).Value = 2073.4
$ws.Cells.Item(99, 9).Value = 1708.5
$ws.Cells.Item(99, 10).Value = 2316.6667
$ws.Cells.Item(99, 11).Value = 1708.5
$ws.Cells.Item(99, 12).Value = 2316.6667
$ws.Cells.Item(99, 13).Value = -210.5
$ws.Cells.Item(99, 14).Value = -5312.6667
$ws.Cells.Item(107, 8).Value = 1446.75
$ws.Cells.Item(107, 9).Value = 1906.9474
$ws.Cells.Item(107, 10).Value = 475.22223
$ws.Cells.Item(107, 11).Value = 1906.9474
$ws.Cells.Item(107, 12).Value = 475.22223
$ws.Cells.Item(107, 13).Value = 13.05259999999998
$ws.Cells.Item(107, 14).Value = -4315.22223
$ws.Cells.Item(122, 8).Value = 4644.24
$ws.Cells.Item(122, 9).Value = 6132
$ws.Cells.Item(122, 11).Value = 18396
$ws.Cells.Item(122, 13).Value = -15946
$ws.Cells.Item(126, 8).Value = 2073.4
$ws.Cells.Item(126, 9).Value = 1708.5
$ws.Cells.Item(126, 10).Value = 2316.6667
$ws.Cells.Item(126, 11).Value = 5125.5
$ws.Cells.Item(126, 12).Value = 6950.000100000001
$ws.Cells.Item(126, 13).Value = -2655.5
$ws.Cells.Item(126, 14).Value = -11890.0001

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(131, 8).Value = 966.7245
$ws.Cells.Item(131, 9).Value = 800
$ws.Cells.Item(131, 10).Value = 971.9895
$ws.Cells.Item(131, 11).Value = 2400
$ws.Cells.Item(131, 12).Value = 2915.9685
$ws.Cells.Item(131, 13).Value = 2640
$ws.Cells.Item(131, 14).Value = -12995.9685

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102, 8).Value = 1169.5454
$ws.Cells.Item(102, 9).Value = 1062.4
$ws.Cells.Item(102, 10).Value = 1399.1428
$ws.Cells.Item(102, 11).Value = 1062.4
$ws.Cells.Item(102, 12).Value = 1399.1428
$ws.Cells.Item(102, 13).Value = 559.5999999999999
$ws.Cells.Item(102, 14).Value = -4643.1428
$ws.Cells.Item(122, 8).Value = 2828.92
$ws.Cells.Item(122, 9).Value = 795.7778
$ws.Cells.Item(122, 10).Value = 3972.5625
$ws.Cells.Item(122, 11).Value = 2387.3334
$ws.Cells.Item(122, 12).Value = 11917.6875
$ws.Cells.Item(122, 13).Value = 62.66660000000002
$ws.Cells.Item(122, 14).Value = -16817.6875
$ws.Cells.Item(126, 8).Value = 1888.7742
$ws.Cells.Item(126, 9).Value = 1877.9
$ws.Cells.Item(126, 10).Value = 1908.5454
$ws.Cells.Item(126, 11).Value = 5633.700000000001
$ws.Cells.Item(126, 12).Value = 5725.6362
$ws.Cells.Item(126, 13).Value = -3163.700000000001
$ws.Cells.Item(126, 14).Value = -10665.6362

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 2238.9375
$ws.Cells.Item(40, 9).Value = 1687.3
$ws.Cells.Item(40, 10).Value = 3158.3333
$ws.Cells.Item(40, 11).Value = 1687.3
$ws.Cells.Item(40, 12).Value = 3158.3333
$ws.Cells.Item(40, 13).Value = -1551.3
$ws.Cells.Item(40, 14).Value = -3430.3333
$ws.Cells.Item(93, 8).Value = 1568.1818
$ws.Cells.Item(93, 9).Value = 1275.5
$ws.Cells.Item(93, 10).Value = 1919.4
$ws.Cells.Item(93, 11).Value = 1275.5
$ws.Cells.Item(93, 12).Value = 1919.4
$ws.Cells.Item(93, 13).Value = -27.5
$ws.Cells.Item(93, 14).Value = -4415.4
$ws.Cells.Item(97, 8).Value = 14696
$ws.Cells.Item(97, 10).Value = 14696
$ws.Cells.Item(97, 12).Value = 14696
$ws.Cells.Item(97, 14).Value = -16678
$ws.Cells.Item(136, 8).Value = 3034.238
$ws.Cells.Item(136, 9).Value = 2445.1304
$ws.Cells.Item(136, 10).Value = 3747.3684
$ws.Cells.Item(136, 11).Value = 7335.3912
$ws.Cells.Item(136, 12).Value = 11242.1052
$ws.Cells.Item(136, 13).Value = -4785.3912
$ws.Cells.Item(136, 14).Value = -16342.1052

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(92, 8).Value = 35550
$ws.Cells.Item(92, 10).Value = 35550
$ws.Cells.Item(92, 12).Value = 35550
$ws.Cells.Item(92, 14).Value = -40542
